# edit.ps1 - apply the "Processing Time & Parrallelization" slide insertion
# and the footer date bump (8/28/2019 -> 8/29/2019) to the HNSCC deck.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "Date Placeholder" text bump: 8/28/2019 -> 8/29/2019
#    This placeholder is cached on the slide master and on every slide
#    layout (not on the slides themselves), so walk all of them.
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "8/29/2019"
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Insert the new "Processing Time & Parrallelization" slide right after
#    the "Probit Fitting Convergence" slide (position 6), using the same
#    "Title and Content" layout as its neighbours.
# ---------------------------------------------------------------------------
$newSlide = $p.Slides.Add(6, 2)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Processing Time & Parrallelization"

$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "`r[20/21] - Total time elapsed: 0:25:10.298157 `r~75 sec / panel `r`rProcessing can parallelize to panel granularity `rassuming available_threads > num_panels, parallelization will reduce time complexity to O(c)`rOtherwise, O(n) `r"

$body.Paragraphs(3,1).IndentLevel = 3
$body.Paragraphs(4,1).IndentLevel = 3
$body.Paragraphs(6,1).IndentLevel = 2
$body.Paragraphs(7,1).IndentLevel = 2
$body.Paragraphs(8,1).IndentLevel = 2
